$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting so numeric-looking / percentage-looking strings
# (e.g. "330.41", "-0.47%") are stored as text, matching the source data
# (all cells in this sheet are inline strings, not numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '330.41'
$ws.Range('E2').Value = '-0.47%'
$ws.Range('D3').Value = '41.42'
$ws.Range('E3').Value = '0.23%'
$ws.Range('D4').Value = '5.631'
$ws.Range('E4').Value = '-1.41%'
$ws.Range('D5').Value = '0.08335'
$ws.Range('E5').Value = '2.73%'
$ws.Range('D6').Value = '8.781'
$ws.Range('E6').Value = '0.45%'
$ws.Range('D7').Value = '1.986'
$ws.Range('E7').Value = '-3.61%'
$ws.Range('D8').Value = '4.479'
$ws.Range('E8').Value = '-1.02%'
$ws.Range('E9').Value = '-1.71%'
$ws.Range('D10').Value = '0.9256'
$ws.Range('E10').Value = '0.11%'
$ws.Range('D11').Value = '0.1294'
$ws.Range('E11').Value = '2.29%'
$ws.Range('D12').Value = '0.1958'
$ws.Range('E12').Value = '-0.09%'
$ws.Range('D13').Value = '0.09343'
$ws.Range('E13').Value = '1.29%'
$ws.Range('D14').Value = '0.03911'
$ws.Range('E14').Value = '4.57%'
$ws.Range('D15').Value = '0.1059'
$ws.Range('E15').Value = '0.59%'
$ws.Range('D16').Value = '0.001303'
$ws.Range('E16').Value = '-0.57%'
$ws.Range('D17').Value = '0.006091'
$ws.Range('E17').Value = '-0.68%'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').Value = '3.441'
$ws.Range('E18').Value = '1.80%'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').Value = '0.3536'
$ws.Range('E19').Value = '0.15%'
$ws.Range('B20').Value = 'MCDex'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D20').Value = '8.546'
$ws.Range('E20').Value = '-2.99%'
$ws.Range('B21').Value = 'ProBitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D21').Value = '0.1370'
$ws.Range('E21').Value = '-3.33%'
$ws.Range('B22').Value = 'ZBToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D22').Value = '0.2420'
$ws.Range('E22').Value = '-7.28%'
$ws.Range('B23').Value = 'CoinExToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D23').Value = '0.04400'
$ws.Range('E23').Value = '-0.74%'
$ws.Range('B24').Value = 'BitKan'
$ws.Range('C24').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D24').Value = '0.001275'
$ws.Range('E24').Value = '1.20%'
$ws.Range('B25').Value = 'HotbitToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D25').Value = '0.004387'
$ws.Range('E25').Value = '-1.42%'
$ws.Range('D26').Value = '0.0001201'
$ws.Range('E26').Value = '-3.18%'
$ws.Range('D39').Value = '0.02814'
$ws.Range('E39').Value = '-0.20%'
$ws.Range('D40').Value = '0.05507'
$ws.Range('E40').Value = '-0.09%'
$ws.Range('D41').Value = '0.007945'
$ws.Range('E41').Value = '3.34%'
$ws.Range('D42').Value = '0.1435'
$ws.Range('E42').Value = '0.97%'
$ws.Range('D43').Value = '0.009321'
$ws.Range('E43').Value = '-6.01%'
$ws.Range('D44').Value = '0.002141'
$ws.Range('E44').Value = '2.44%'
$ws.Range('D45').Value = '0.01108'
$ws.Range('E45').Value = '0.88%'
$ws.Range('D46').Value = '0.00007092'
$ws.Range('E46').Value = '4.43%'
$ws.Range('E47').Value = '0.04%'
$ws.Range('D48').Value = '0.003429'
$ws.Range('E48').Value = '14.58%'
$ws.Range('E49').Value = '-0.02%'
$ws.Range('D50').Value = '0.00002101'
$ws.Range('E50').Value = '0.04%'
$ws.Range('D51').Value = '0.0002001'
$ws.Range('E51').Value = '0.04%'
